$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 9.016222821604932
$ws.Range("C2").Value = 4.566638621406133
$ws.Range("D2").Value = 5.990388079305196
$ws.Range("E2").Value = 12.24725207756675
$ws.Range("G2").Value = 3.644644103155783
$ws.Range("I2").Value = 21.94064994617839
$ws.Range("K2").Value = 9.178072607331282
$ws.Range("M2").Value = 13.94576580312774
$ws.Range("N2").Value = 19.06162539099923
$ws.Range("O2").Value = 23.75267136655931
$ws.Range("B3").Value = 8.745277880832385
$ws.Range("C3").Value = 4.304111141652966
$ws.Range("D3").Value = 5.872487292434665
$ws.Range("E3").Value = 12.03044354208856
$ws.Range("G3").Value = 3.646594319692972
$ws.Range("I3").Value = 22.0057065663292
$ws.Range("K3").Value = 8.998737671847719
$ws.Range("M3").Value = 13.78027426953461
$ws.Range("N3").Value = 19.1231232535524
$ws.Range("O3").Value = 23.79536246083149
$ws.Range("B4").Value = 8.576394002834267
$ws.Range("C4").Value = 4.133599557423531
$ws.Range("D4").Value = 5.800648373065957
$ws.Range("E4").Value = 11.89955102371895
$ws.Range("G4").Value = 3.647855411276372
$ws.Range("I4").Value = 22.0498093504767
$ws.Range("K4").Value = 8.888649174604984
$ws.Range("M4").Value = 13.6811957368081
$ws.Range("N4").Value = 19.16263027857524
$ws.Range("O4").Value = 23.8267152334998
$ws.Range("B5").Value = 8.507051097611948
$ws.Range("C5").Value = 4.061789297869836
$ws.Range("D5").Value = 5.771557700178214
$ws.Range("E5").Value = 11.84684985933928
$ws.Range("G5").Value = 3.648385372958467
$ws.Range("I5").Value = 22.06882531613813
$ws.Range("K5").Value = 8.843856742378044
$ws.Range("M5").Value = 13.64150437951378
$ws.Range("N5").Value = 19.17917030099943
$ws.Range("O5").Value = 23.84078169233001
$ws.Range("B6").Value = 8.495508719406706
$ws.Range("C6").Value = 4.049725431546673
$ws.Range("D6").Value = 5.766739677208825
$ws.Range("E6").Value = 11.83813970844915
$ws.Range("G6").Value = 3.648474343977018
$ws.Range("I6").Value = 22.07204588791978
$ws.Range("K6").Value = 8.836425070893778
$ws.Range("M6").Value = 13.63495630764745
$ws.Range("N6").Value = 19.18194340695739
$ws.Range("O6").Value = 23.84319524426875
$ws.Range("B7").Value = 8.57546077636211
$ws.Range("C7").Value = 4.13264048739919
$ws.Range("D7").Value = 5.800255239805523
$ws.Range("E7").Value = 11.89883758905836
$ws.Range("G7").Value = 3.647862493444218
$ws.Range("I7").Value = 22.05006158286194
$ws.Range("K7").Value = 8.888044721926642
$ws.Range("M7").Value = 13.68065761526152
$ws.Range("N7").Value = 19.16285155750248
$ws.Range("O7").Value = 23.82689971944568
$ws.Range("B8").Value = 8.923390008056119
$ws.Range("C8").Value = 4.47806244718128
$ws.Range("D8").Value = 5.949648211373358
$ws.Range("E8").Value = 12.17208098387965
$ws.Range("G8").Value = 3.645303357915646
$ws.Range("I8").Value = 21.96221747475769
$ws.Range("K8").Value = 9.11627067966875
$ws.Range("M8").Value = 13.88820687076865
$ws.Range("N8").Value = 19.08246801357188
$ws.Range("O8").Value = 23.76632291308441
$ws.Range("B9").Value = 9.581090733146254
$ws.Range("C9").Value = 5.080941221290415
$ws.Range("D9").Value = 6.245034905400384
$ws.Range("E9").Value = 12.72214505987117
$ws.Range("G9").Value = 3.640787607347804
$ws.Range("I9").Value = 21.82302314977187
$ws.Range("K9").Value = 9.561293260030022
$ws.Range("M9").Value = 14.3132021480202
$ws.Range("N9").Value = 18.93863942794092
$ws.Range("O9").Value = 23.68842282467934
$ws.Range("B10").Value = 10.04376085559096
$ws.Range("C10").Value = 5.47798125044414
$ws.Range("D10").Value = 6.460991099861824
$ws.Range("E10").Value = 13.13037401220226
$ws.Range("G10").Value = 3.637773081971221
$ws.Range("I10").Value = 21.74101972543297
$ws.Range("K10").Value = 9.88320896171442
$ws.Range("M10").Value = 14.63352841209302
$ws.Range("N10").Value = 18.8412991249664
$ws.Range("O10").Value = 23.65624390245514
$ws.Range("B11").Value = 10.24877018682809
$ws.Range("C11").Value = 5.648593642689994
$ws.Range("D11").Value = 6.558495961314454
$ws.Range("E11").Value = 13.3160367271202
$ws.Range("G11").Value = 3.636466837628387
$ws.Range("I11").Value = 21.70813500410051
$ws.Range("K11").Value = 10.02786930904975
$ws.Range("M11").Value = 14.78038942005485
$ws.Range("N11").Value = 18.79880740478436
$ws.Range("O11").Value = 23.64706488569925
$ws.Range("B12").Value = 10.32554308778393
$ws.Range("C12").Value = 5.711760137087875
$ws.Range("D12").Value = 6.595273171017493
$ws.Range("E12").Value = 13.38626111341368
$ws.Range("G12").Value = 3.635981502383168
$ws.Range("I12").Value = 21.69631948157615
$ws.Range("K12").Value = 10.08233852280263
$ws.Range("M12").Value = 14.83611483379398
$ws.Range("N12").Value = 18.78297281183098
$ws.Range("O12").Value = 23.64437499433103
$ws.Range("B13").Value = 10.30904800684297
$ws.Range("C13").Value = 5.698220198008977
$ws.Range("D13").Value = 6.587359639048771
$ws.Range("E13").Value = 13.37114185998009
$ws.Range("G13").Value = 3.63608561464174
$ws.Range("I13").Value = 21.69883579877798
$ws.Range("K13").Value = 10.07062224685958
$ws.Range("M13").Value = 14.82410919271171
$ws.Range("N13").Value = 18.78637170265828
$ws.Range("O13").Value = 23.64491934197595
$ws.Range("B14").Value = 10.25510396746778
$ws.Range("C14").Value = 5.653819269692606
$ws.Range("D14").Value = 6.56152478777334
$ws.Range("E14").Value = 13.32181619781234
$ws.Range("G14").Value = 3.636426722465248
$ws.Range("I14").Value = 21.70715015362649
$ws.Range("K14").Value = 10.03235704799351
$ws.Range("M14").Value = 14.78497202722445
$ws.Range("N14").Value = 18.79749955716191
$ws.Range("O14").Value = 23.64682783001804
$ws.Range("B15").Value = 10.22194765457869
$ws.Range("C15").Value = 5.626434755036536
$ws.Range("D15").Value = 6.545680059470758
$ws.Range("E15").Value = 13.29158991980606
$ws.Range("G15").Value = 3.636636872050703
$ws.Range("I15").Value = 21.7123259698558
$ws.Range("K15").Value = 10.00887650769863
$ws.Range("M15").Value = 14.76101248185417
$ws.Range("N15").Value = 18.80434900744981
$ws.Range("O15").Value = 23.64809921628075
$ws.Range("B16").Value = 10.03024658424573
$ws.Range("C16").Value = 5.466629808700444
$ws.Range("D16").Value = 6.45460066473445
$ws.Range("E16").Value = 13.11823309863124
$ws.Range("G16").Value = 3.637859753623315
$ws.Range("I16").Value = 21.74325789378273
$ws.Range("K16").Value = 9.873714508489588
$ws.Range("M16").Value = 14.62394923635498
$ws.Range("N16").Value = 18.84411195792141
$ws.Range("O16").Value = 23.65695371971086
$ws.Range("B17").Value = 9.911190011106203
$ws.Range("C17").Value = 5.366030458075886
$ws.Range("D17").Value = 6.398509001873308
$ws.Range("E17").Value = 13.01181839834663
$ws.Range("G17").Value = 3.638626586394861
$ws.Range("I17").Value = 21.76336672697361
$ws.Range("K17").Value = 9.79030117818456
$ws.Range("M17").Value = 14.54012089241629
$ws.Range("N17").Value = 18.86896259962078
$ws.Range("O17").Value = 23.66378474780913
$ws.Range("B18").Value = 9.842201686354249
$ws.Range("C18").Value = 5.307226851721507
$ws.Range("D18").Value = 6.366179674193624
$ws.Range("E18").Value = 12.95061182165826
$ws.Range("G18").Value = 3.639073776596916
$ws.Range("I18").Value = 21.77534863805457
$ws.Range("K18").Value = 9.742159489544497
$ws.Range("M18").Value = 14.49201570772338
$ws.Range("N18").Value = 18.88342447581522
$ws.Range("O18").Value = 23.66822760128165
$ws.Range("B19").Value = 9.818758199572668
$ws.Range("C19").Value = 5.287155350650352
$ws.Range("D19").Value = 6.355223289476614
$ws.Range("E19").Value = 12.92989074864543
$ws.Range("G19").Value = 3.639226241492555
$ws.Range("I19").Value = 21.77947687370946
$ws.Range("K19").Value = 9.725832990391048
$ws.Range("M19").Value = 14.4757487214236
$ws.Range("N19").Value = 18.88834998397759
$ws.Range("O19").Value = 23.66982008958657
$ws.Range("B20").Value = 9.923917145983484
$ws.Range("C20").Value = 5.376836896184212
$ws.Range("D20").Value = 6.404487275235453
$ws.Range("E20").Value = 13.02314694895569
$ws.Range("G20").Value = 3.638544321825706
$ws.Range("I20").Value = 21.7611830562487
$ws.Range("K20").Value = 9.799198104057503
$ws.Range("M20").Value = 14.54903348880847
$ws.Range("N20").Value = 18.866299781293
$ws.Range("O20").Value = 23.6630043858067
$ws.Range("B21").Value = 10.27097252542588
$ws.Range("C21").Value = 5.666900009516747
$ws.Range("D21").Value = 6.569117374623025
$ws.Range("E21").Value = 13.33630715684108
$ws.Range("G21").Value = 3.636326278500615
$ws.Range("I21").Value = 21.70469071888652
$ws.Range("K21").Value = 10.0436053101743
$ws.Range("M21").Value = 14.79646492341201
$ws.Range("N21").Value = 18.79422409510493
$ws.Range("O21").Value = 23.64624592331871
$ws.Range("B22").Value = 10.49275321208831
$ws.Range("C22").Value = 5.848075434588329
$ws.Range("D22").Value = 6.675846006069107
$ws.Range("E22").Value = 13.54046426055145
$ws.Range("G22").Value = 3.634930909826749
$ws.Range("I22").Value = 21.67148459571616
$ws.Range("K22").Value = 10.20150659599365
$ws.Range("M22").Value = 14.95880693050632
$ws.Range("N22").Value = 18.74861069007263
$ws.Range("O22").Value = 23.63987488018763
$ws.Range("B23").Value = 10.37486844989643
$ws.Range("C23").Value = 5.752147434905095
$ws.Range("D23").Value = 6.618974902872414
$ws.Range("E23").Value = 13.43157279806092
$ws.Range("G23").Value = 3.635670695842984
$ws.Range("I23").Value = 21.6888668598923
$ws.Range("K23").Value = 10.11741659045181
$ws.Range("M23").Value = 14.87212115996489
$ws.Range("N23").Value = 18.77281925752399
$ws.Range("O23").Value = 23.64285580076466
$ws.Range("B24").Value = 9.918164888024375
$ws.Range("C24").Value = 5.371954319540291
$ws.Range("D24").Value = 6.401784748592156
$ws.Range("E24").Value = 13.01802539060749
$ws.Range("G24").Value = 3.638581493922965
$ws.Range("I24").Value = 21.76216898226012
$ws.Range("K24").Value = 9.795176379821674
$ws.Range("M24").Value = 14.54500382325053
$ws.Range("N24").Value = 18.8675030965052
$ws.Range("O24").Value = 23.66335558143097
$ws.Range("B25").Value = 9.406402757674833
$ws.Range("C25").Value = 4.925882898566395
$ws.Range("D25").Value = 6.165133342172064
$ws.Range("E25").Value = 12.57229208024554
$ws.Range("G25").Value = 3.641955758387346
$ws.Range("I25").Value = 21.85712779215719
$ws.Range("K25").Value = 9.441544046636087
$ws.Range("M25").Value = 14.19659983887774
$ws.Range("N25").Value = 18.97607988836647
$ws.Range("O25").Value = 23.70510475889124
